$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title shape: "A" + " " + "slide" -> single run "A slide"
$titleShape = $s.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Text = "__tmp__"
$titleRange.Text = "A slide"

# Table cell: "a" + " " + "table" -> single run "a table"
$tableShape = $s.Shapes.Item(3)
$tbl = $tableShape.Table
$cell = $tbl.Cell(1, 2)
$cellRange = $cell.Shape.TextFrame.TextRange
$cellRange.Text = "__tmp__"
$cellRange.Text = "a table"

# TextBox: "Plus" + " " + "an" + " " + "image" -> single run "Plus an image"
$textboxShape = $s.Shapes.Item(7)
$textboxRange = $textboxShape.TextFrame.TextRange
$textboxRange.Text = "__tmp__"
$textboxRange.Text = "Plus an image"
